$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temp data added to rows 2-7 for columns A (Date) through F (retweets).
$data = @(
    @(43922, 1, 1, 2, 4, 2),
    @(43923, 1, 4, 3, 6, 3),
    @(43924, 1, 5, 3, 6, 2),
    @(43925, 1, 3, 2, 6, 3),
    @(43926, 1, 3, 1, 6, 2),
    @(43927, 1, 2, 1, 6, 3)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Count; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $values[$j]
    }
}

# Update the active selection to match the author's last selected cell.
$ws.Range("G14").Select() | Out-Null
